$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("H4").Value = 83
$ws.Range("J4").Value = 100
$ws.Range("L4").Value = 100
$ws.Range("N4").Value = -328
$ws.Range("H9").Value = 65.5
$ws.Range("I9").Value = 65
$ws.Range("K9").Value = 65
$ws.Range("M9").Value = 104
$ws.Range("H11").Value = 40.333332
$ws.Range("I11").Value = 40.333332
$ws.Range("K11").Value = 40.333332
$ws.Range("M11").Value = 99.666668
$ws.Range("H33").Value = 104.8
$ws.Range("I33").Value = 74.77778000000001
$ws.Range("K33").Value = 74.77778000000001
$ws.Range("M33").Value = 154.22222
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").Value = ""
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").Value = ""
$ws.Range("H70").Value = 3884
$ws.Range("I70").Value = 3300
$ws.Range("J70").Value = 4662.6665
$ws.Range("K70").Value = 9900
$ws.Range("L70").Value = 13987.9995
$ws.Range("M70").Value = -9630
$ws.Range("N70").Value = -14527.9995
$ws.Range("H73").Value = 3884
$ws.Range("I73").Value = 3300
$ws.Range("J73").Value = 4662.6665
$ws.Range("K73").Value = 9900
$ws.Range("L73").Value = 13987.9995
$ws.Range("M73").Value = -8964
$ws.Range("N73").Value = -15859.9995
$ws.Range("H80").Value = 1579.8572
$ws.Range("J80").Value = 1758.3334
$ws.Range("L80").Value = 5275.0002
$ws.Range("N80").Value = -7271.0002
$ws.Range("H83").Value = 1579.8572
$ws.Range("J83").Value = 1758.3334
$ws.Range("L83").Value = 15825.0006
$ws.Range("N83").Value = -25809.0006
$ws.Range("H88").Value = 6712.143
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 6712.143
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 6712.143
$ws.Range("M88").Value = ""
$ws.Range("N88").Value = -7524.143
$ws.Range("H91").Value = 6712.143
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 6712.143
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 6712.143
$ws.Range("M91").Value = ""
$ws.Range("N91").Value = -9520.143
$ws.Range("H100").Value = 801.6667
$ws.Range("I100").Value = 801.6667
$ws.Range("K100").Value = 801.6667
$ws.Range("M100").Value = -260.6667
$ws.Range("H101").Value = 5666.6665
$ws.Range("H138").Value = 3114.5417
$ws.Range("I138").Value = 2721.1428
$ws.Range("J138").Value = 3276.5293
$ws.Range("K138").Value = 8163.428400000001
$ws.Range("L138").Value = 9829.5879
$ws.Range("M138").Value = -3023.428400000001
$ws.Range("N138").Value = -20109.5879
$ws = $wb.Worksheets.Item(2)
$ws.Range("H15").Value = 3949.5
$ws.Range("J15").Value = 3949.5
$ws.Range("L15").Value = 3949.5
$ws.Range("N15").Value = -4649.5
$ws.Range("H45").Value = 4002.1428
$ws.Range("I45").Value = 3856.6667
$ws.Range("J45").Value = 4875
$ws.Range("K45").Value = 3856.6667
$ws.Range("L45").Value = 4875
$ws.Range("M45").Value = -3479.6667
$ws.Range("N45").Value = -5629
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").Value = ""
$ws.Range("H103").Value = 29990
$ws.Range("J103").Value = 29990
$ws.Range("L103").Value = 29990
$ws.Range("N103").Value = -32334
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").Value = ""
$ws.Range("H122").Value = 1012
$ws.Range("I122").Value = 1012
$ws.Range("K122").Value = 3036
$ws.Range("M122").Value = -586
$ws = $wb.Worksheets.Item(3)
$ws.Range("H33").Value = 950
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").Value = ""
$ws.Range("H75").Value = 8000
$ws.Range("I75").Value = 8000
$ws.Range("K75").Value = 8000
$ws.Range("M75").Value = -7064
$ws.Range("H78").Value = 8000
$ws.Range("I78").Value = 8000
$ws.Range("K78").Value = 24000
$ws.Range("M78").Value = -19320
$ws.Range("H134").Value = 4888.6665
$ws = $wb.Worksheets.Item(4)
$ws.Range("H31").Value = 9611.111000000001
$ws.Range("I31").Value = 8250
$ws.Range("K31").Value = 8250
$ws.Range("M31").Value = -7955
$ws.Range("H32").Value = 2570.1538
$ws.Range("I32").Value = 789.2
$ws.Range("J32").Value = 3683.25
$ws.Range("K32").Value = 789.2
$ws.Range("L32").Value = 3683.25
$ws.Range("M32").Value = -473.2
$ws.Range("N32").Value = -4315.25
$ws.Range("H34").Value = 9611.111000000001
$ws.Range("I34").Value = 8250
$ws.Range("K34").Value = 8250
$ws.Range("M34").Value = -8048
$ws.Range("H59").Value = 30339
$ws.Range("J59").Value = 37900
$ws.Range("L59").Value = 37900
$ws.Range("N59").Value = -40190
$ws.Range("H60").Value = 2950
$ws.Range("I60").Value = 2950
$ws.Range("K60").Value = 2950
$ws.Range("M60").Value = -2439
$ws.Range("H68").Value = 44500
$ws.Range("I68").Value = 40000
$ws.Range("J68").Value = 49000
$ws.Range("K68").Value = 40000
$ws.Range("L68").Value = 49000
$ws.Range("M68").Value = -39251
$ws.Range("N68").Value = -50498
$ws.Range("H71").Value = 44500
$ws.Range("I71").Value = 40000
$ws.Range("J71").Value = 49000
$ws.Range("K71").Value = 120000
$ws.Range("L71").Value = 147000
$ws.Range("M71").Value = -116256
$ws.Range("N71").Value = -154488
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").Value = ""
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").Value = ""
$ws.Range("H105").Value = 3332
$ws.Range("I105").Value = 3332
$ws.Range("K105").Value = 3332
$ws.Range("M105").Value = -1585
$ws.Range("H117").Value = 0
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("M117").Value = ""
$ws.Range("N117").Value = ""
$ws.Range("H134").Value = 8825
$ws.Range("I134").Value = 650
$ws.Range("K134").Value = 1950
$ws.Range("M134").Value = 585
$ws = $wb.Worksheets.Item(5)
$ws.Range("H12").Value = 30.285715
$ws.Range("J12").Value = 24.1
$ws.Range("L12").Value = 72.30000000000001
$ws.Range("N12").Value = -418.3
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").Value = ""
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").Value = ""
$ws.Range("H80").Value = 1003
$ws.Range("J80").Value = 1003
$ws.Range("L80").Value = 3009
$ws.Range("N80").Value = -4881
$ws.Range("H83").Value = 1003
$ws.Range("J83").Value = 1003
$ws.Range("L83").Value = 9027
$ws.Range("N83").Value = -18387
$ws.Range("H121").Value = 1582.5
$ws.Range("J121").Value = 1582.5
$ws.Range("L121").Value = 4747.5
$ws.Range("N121").Value = -7367.5
$ws = $wb.Worksheets.Item(6)
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").Value = ""
$ws.Range("H126").Value = 2000
$ws.Range("I126").Value = 2000
$ws.Range("K126").Value = 6000
$ws.Range("M126").Value = -3530
$ws = $wb.Worksheets.Item(7)
$ws.Range("H7").Value = 2998.75
$ws.Range("I7").Value = 2000
$ws.Range("J7").Value = 3331.6667
$ws.Range("K7").Value = 2000
$ws.Range("L7").Value = 3331.6667
$ws.Range("M7").Value = -1888
$ws.Range("N7").Value = -3555.6667
$ws.Range("H16").Value = 899.5
$ws.Range("I16").Value = 899.5
$ws.Range("K16").Value = 899.5
$ws.Range("M16").Value = -729.5
$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").Value = ""
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = ""
$ws.Range("N40").Value = ""
$ws.Range("H46").Value = 4543.222
$ws.Range("I46").Value = 4948.4287
$ws.Range("K46").Value = 4948.4287
$ws.Range("M46").Value = -4760.4287
$ws.Range("H93").Value = 1145.4546
$ws.Range("I93").Value = 1220
$ws.Range("J93").Value = 400
$ws.Range("K93").Value = 1220
$ws.Range("L93").Value = 400
$ws.Range("M93").Value = 28
$ws.Range("N93").Value = -2896
$ws.Range("H126").Value = 2998.75
$ws.Range("I126").Value = 2000
$ws.Range("J126").Value = 3331.6667
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 9995.000100000001
$ws.Range("M126").Value = -3530
$ws.Range("N126").Value = -14935.0001
$ws = $wb.Worksheets.Item(8)
$ws.Range("H45").Value = 33650
$ws.Range("I45").Value = 32975
$ws.Range("K45").Value = 32975
$ws.Range("M45").Value = -32484
$ws.Range("H96").Value = 2999
$ws.Range("I96").Value = 2999
$ws.Range("K96").Value = 2999
$ws.Range("M96").Value = -1626
$ws.Range("H122").Value = 4061
$ws.Range("I122").Value = 4225
$ws.Range("K122").Value = 12675
$ws.Range("M122").Value = -10225
